# updating with august charges
$wb = $excel.ActiveWorkbook

# --- eto_use_alt (sheet3): append August charge rows 51-54 ---
$wsAlt = $wb.Worksheets.Item("eto_use_alt")

$wsAlt.Range("A50").Copy() | Out-Null
$wsAlt.Range("A51").PasteSpecial(-4122) | Out-Null
$wsAlt.Range("A51").Value = 45874
$wsAlt.Range("B51").Value = "CL007"

$wsAlt.Range("A50").Copy() | Out-Null
$wsAlt.Range("A52").PasteSpecial(-4122) | Out-Null
$wsAlt.Range("A52").Value = 45876
$wsAlt.Range("B52").Value = "CL015"

$wsAlt.Range("A50").Copy() | Out-Null
$wsAlt.Range("A53").PasteSpecial(-4122) | Out-Null
$wsAlt.Range("A53").Value = 45883
$wsAlt.Range("B53").Value = "CL015"

$wsAlt.Range("A50").Copy() | Out-Null
$wsAlt.Range("A54").PasteSpecial(-4122) | Out-Null
$wsAlt.Range("A54").Value = 45890
$wsAlt.Range("B54").Value = "CL015"

# --- fluoro_use (sheet4): fill row 6, append rows 7-9 ---
$wsFluoro = $wb.Worksheets.Item("fluoro_use")

$wsFluoro.Range("A6").Value = 45881
$wsFluoro.Range("B6").Value = "CL014"

$wsFluoro.Range("A5").Copy() | Out-Null
$wsFluoro.Range("A7").PasteSpecial(-4122) | Out-Null
$wsFluoro.Range("A7").Value = 45887
$wsFluoro.Range("B7").Value = "CL014"

$wsFluoro.Range("A5").Copy() | Out-Null
$wsFluoro.Range("A8").PasteSpecial(-4122) | Out-Null
$wsFluoro.Range("A8").Value = 45888
$wsFluoro.Range("B8").Value = "CL014"

$wsFluoro.Range("A5").Copy() | Out-Null
$wsFluoro.Range("A9").PasteSpecial(-4122) | Out-Null
$wsFluoro.Range("A9").Value = 45894
$wsFluoro.Range("B9").Value = "CL014"

# --- view state: fluoro_use becomes the active/selected tab ---
$wsFluoro.Activate()
$wsAlt.Range("B54").Select() | Out-Null
$wsFluoro.Activate()
$wsFluoro.Range("B4:B9").Select() | Out-Null

Write-Host "done"
